$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 "time_taken" - copy formatting from neighboring header E1 (panel)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data rows: plain inline-string values, no special style
$ws.Range("F2").Value = "2021-10-05 13:42:04.318923"
$ws.Range("F3").Value = "2021-10-05 13:42:04.318935"
